$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 544900
$ws.Range("L3").Value = 44900

$ws.Range("K4").Value = 107100.3
$ws.Range("L4").Value = 6200.3
